$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (General/numeric-looking) Price values to stay as text,
# matching the source data which stores these as strings, not numbers.
$forceTextCells = @("D4","D5","D7","D8","D9","D10","D11","D13","D14","D15","D16","D18","D20","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D47","D48","D49","D50","D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "27.706.57"
$ws.Range("E2").Value = "  -0.37%  "

# Row 3
$ws.Range("D3").Value = "1.850.33"
$ws.Range("E3").Value = "  -0.88%  "

# Row 4
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  -2.44%  "

# Row 5
$ws.Range("D5").Value = "319.26"
$ws.Range("E5").Value = "  -1.61%  "

# Row 6
$ws.Range("E6").Value = "  -2.33%  "

# Row 7
$ws.Range("D7").Value = "0.4334"
$ws.Range("E7").Value = "  -2.13%  "

# Row 8
$ws.Range("D8").Value = "0.3768"
$ws.Range("E8").Value = "  -0.82%  "

# Row 9
$ws.Range("D9").Value = "0.07382"
$ws.Range("E9").Value = "  -1.13%  "

# Row 10
$ws.Range("D10").Value = "0.8835"
$ws.Range("E10").Value = "  -0.37%  "

# Row 11
$ws.Range("D11").Value = "21.65"
$ws.Range("E11").Value = "  -0.67%  "

# Row 12
$ws.Range("D12").Value = "1.854.90"
$ws.Range("E12").Value = "  -0.95%  "

# Row 13
$ws.Range("D13").Value = "6.751"
$ws.Range("E13").Value = "  -0.21%  "

# Row 14
$ws.Range("D14").Value = "5.475"
$ws.Range("E14").Value = "  -1.59%  "

# Row 15
$ws.Range("D15").Value = "0.07152"
$ws.Range("E15").Value = "  -1.16%  "

# Row 16
$ws.Range("D16").Value = "88.15"
$ws.Range("E16").Value = "  +5.08%  "

# Row 17
$ws.Range("E17").Value = "  -2.41%  "

# Row 18
$ws.Range("D18").Value = "0.000009033"
$ws.Range("E18").Value = "  -1.34%  "

# Row 19
$ws.Range("E19").Value = "  -2.34%  "

# Row 20
$ws.Range("D20").Value = "15.53"
$ws.Range("E20").Value = "  -0.33%  "

# Row 21
$ws.Range("D21").Value = "27.714.46"
$ws.Range("E21").Value = "  -0.45%  "

# Row 22
$ws.Range("D22").Value = "5.265"
$ws.Range("E22").Value = "  -1.20%  "

# Row 23
$ws.Range("D23").Value = "11.22"
$ws.Range("E23").Value = "  -1.40%  "

# Row 24
$ws.Range("D24").Value = "2.083.56"
$ws.Range("E24").Value = "  -0.64%  "

# Row 25
$ws.Range("D25").Value = "2.030"
$ws.Range("E25").Value = "  +2.99%  "

# Row 26
$ws.Range("D26").Value = "155.71"
$ws.Range("E26").Value = "  -1.92%  "

# Row 27
$ws.Range("D27").Value = "18.64"
$ws.Range("E27").Value = "  -1.41%  "

# Row 28
$ws.Range("D28").Value = "2.142"
$ws.Range("E28").Value = "  +7.50%  "

# Row 29
$ws.Range("D29").Value = "5.417"
$ws.Range("E29").Value = "  +1.53%  "

# Row 30
$ws.Range("D30").Value = "120.49"
$ws.Range("E30").Value = "  +2.28%  "

# Row 31
$ws.Range("D31").Value = "0.08957"
$ws.Range("E31").Value = "  -1.64%  "

# Row 32
$ws.Range("D32").Value = "1.236"
$ws.Range("E32").Value = "  +1.32%  "

# Row 33
$ws.Range("D33").Value = "0.7761"
$ws.Range("E33").Value = "  -0.66%  "

# Row 34
$ws.Range("D34").Value = "4.579"
$ws.Range("E34").Value = "  -0.19%  "

# Row 35
$ws.Range("D35").Value = "2.929"
$ws.Range("E35").Value = "  -5.78%  "

# Row 36
$ws.Range("D36").Value = "1.144"
$ws.Range("E36").Value = "  -1.98%  "

# Row 37
$ws.Range("D37").Value = "1.014"
$ws.Range("E37").Value = "  -2.37%  "

# Row 38
$ws.Range("D38").Value = "0.05357"
$ws.Range("E38").Value = "  -0.14%  "

# Row 39
$ws.Range("D39").Value = "0.01973"
$ws.Range("E39").Value = "  -1.51%  "

# Row 40
$ws.Range("D40").Value = "7.153"
$ws.Range("E40").Value = "  +3.75%  "

# Row 41
$ws.Range("D41").Value = "2.869"
$ws.Range("E41").Value = "  +0.43%  "

# Row 42
$ws.Range("D42").Value = "0.5181"
$ws.Range("E42").Value = "  -0.60%  "

# Row 43
$ws.Range("D43").Value = "0.1683"
$ws.Range("E43").Value = "  -0.85%  "

# Row 44
$ws.Range("D44").Value = "8.954"
$ws.Range("E44").Value = "  +2.83%  "

# Row 45
$ws.Range("D45").Value = "110.83"
$ws.Range("E45").Value = "  +0.77%  "

# Row 46
$ws.Range("E46").Value = "  +0.91%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.4748"
$ws.Range("E47").Value = "  +0.55%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.709"
$ws.Range("E48").Value = "  -1.14%  "

# Row 49
$ws.Range("D49").Value = "0.06511"
$ws.Range("E49").Value = "  +0.77%  "

# Row 50
$ws.Range("D50").Value = "1.012"
$ws.Range("E50").Value = "  -2.63%  "

# Row 51
$ws.Range("D51").Value = "1.897"
$ws.Range("E51").Value = "  -0.34%  "
